# Update Leve profit/price figures across all profession sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) to reflect refreshed market
# data pulled by the scheduled runner.

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 978.0769  # H28: was 1017.9167
$ws.Cells.Item(28, 9).Value = 979.2222  # I28: was 1039.125
$ws.Cells.Item(28, 11).Value = 979.2222  # K28: was 1039.125
$ws.Cells.Item(28, 13).Value = -494.2222  # M28: was -554.125
# Row 43
$ws.Cells.Item(43, 8).Value = 5109.6  # H43: was 8332.333000000001
$ws.Cells.Item(43, 9).Value = 700  # I43: was 0
$ws.Cells.Item(43, 10).Value = 6212  # J43: was 8332.333000000001
$ws.Cells.Item(43, 11).Value = 700  # K43: was 0
$ws.Cells.Item(43, 12).Value = 6212  # L43: was 8332.333000000001
$ws.Cells.Item(43, 13).Value = -631  # M43: was None
$ws.Cells.Item(43, 14).Value = -6350  # N43: was -8470.333000000001
# Row 88
$ws.Cells.Item(88, 8).Value = 5025.4  # H88: was 5688.25
$ws.Cells.Item(88, 9).Value = 2374.5  # I88: was 2375
$ws.Cells.Item(88, 11).Value = 2374.5  # K88: was 2375
$ws.Cells.Item(88, 13).Value = -1968.5  # M88: was -1969
# Row 91
$ws.Cells.Item(91, 8).Value = 5025.4  # H91: was 5688.25
$ws.Cells.Item(91, 9).Value = 2374.5  # I91: was 2375
$ws.Cells.Item(91, 11).Value = 2374.5  # K91: was 2375
$ws.Cells.Item(91, 13).Value = -970.5  # M91: was -971
# Row 98
$ws.Cells.Item(98, 8).Value = 652.0833  # H98: was 598.6667
$ws.Cells.Item(98, 9).Value = 438.0909  # I98: was 459.07693
$ws.Cells.Item(98, 10).Value = 3006  # J98: was 1506
$ws.Cells.Item(98, 11).Value = 438.0909  # K98: was 459.07693
$ws.Cells.Item(98, 12).Value = 3006  # L98: was 1506
$ws.Cells.Item(98, 13).Value = 1059.9091  # M98: was 1038.92307
$ws.Cells.Item(98, 14).Value = -6002  # N98: was -4502
# Row 122
$ws.Cells.Item(122, 8).Value = 652.0833  # H122: was 598.6667
$ws.Cells.Item(122, 9).Value = 438.0909  # I122: was 459.07693
$ws.Cells.Item(122, 10).Value = 3006  # J122: was 1506
$ws.Cells.Item(122, 11).Value = 1314.2727  # K122: was 1377.23079
$ws.Cells.Item(122, 12).Value = 9018  # L122: was 4518
$ws.Cells.Item(122, 13).Value = 1135.7273  # M122: was 1072.76921
$ws.Cells.Item(122, 14).Value = -13918  # N122: was -9418
# Row 125
$ws.Cells.Item(125, 8).Value = 3289  # H125: was 3864.75
$ws.Cells.Item(125, 9).Value = 986.5  # I125: was 987
$ws.Cells.Item(125, 11).Value = 8878.5  # K125: was 8883
$ws.Cells.Item(125, 13).Value = -6418.5  # M125: was -6423

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 63
$ws.Cells.Item(63, 8).Value = 4713.5  # H63: was 5265.2144
$ws.Cells.Item(63, 9).Value = 2557.375  # I63: was 2907.7144
$ws.Cells.Item(63, 10).Value = 6869.625  # J63: was 7622.7144
$ws.Cells.Item(63, 11).Value = 2557.375  # K63: was 2907.7144
$ws.Cells.Item(63, 12).Value = 6869.625  # L63: was 7622.7144
$ws.Cells.Item(63, 13).Value = -1871.375  # M63: was -2221.7144
$ws.Cells.Item(63, 14).Value = -8241.625  # N63: was -8994.714400000001
# Row 66
$ws.Cells.Item(66, 8).Value = 4713.5  # H66: was 5265.2144
$ws.Cells.Item(66, 9).Value = 2557.375  # I66: was 2907.7144
$ws.Cells.Item(66, 10).Value = 6869.625  # J66: was 7622.7144
$ws.Cells.Item(66, 11).Value = 12786.875  # K66: was 14538.572
$ws.Cells.Item(66, 12).Value = 34348.125  # L66: was 38113.572
$ws.Cells.Item(66, 13).Value = -9354.875  # M66: was -11106.572
$ws.Cells.Item(66, 14).Value = -41212.125  # N66: was -44977.572
# Row 74
$ws.Cells.Item(74, 8).Value = 1348.3379  # H74: was 1379.096
$ws.Cells.Item(74, 9).Value = 1073.7246  # I74: was 1102.7059
$ws.Cells.Item(74, 11).Value = 1073.7246  # K74: was 1102.7059
$ws.Cells.Item(74, 13).Value = -199.7246  # M74: was -228.7058999999999
# Row 77
$ws.Cells.Item(77, 8).Value = 1348.3379  # H77: was 1379.096
$ws.Cells.Item(77, 9).Value = 1073.7246  # I77: was 1102.7059
$ws.Cells.Item(77, 11).Value = 5368.623  # K77: was 5513.5295
$ws.Cells.Item(77, 13).Value = -1000.623  # M77: was -1145.5295
# Row 97
$ws.Cells.Item(97, 8).Value = 1118  # H97: was 1883.9
$ws.Cells.Item(97, 9).Value = 663.8889  # I97: was 708.375
$ws.Cells.Item(97, 10).Value = 3161.5  # J97: was 6586
$ws.Cells.Item(97, 11).Value = 663.8889  # K97: was 708.375
$ws.Cells.Item(97, 12).Value = 3161.5  # L97: was 6586
$ws.Cells.Item(97, 13).Value = -167.8889  # M97: was -212.375
$ws.Cells.Item(97, 14).Value = -4153.5  # N97: was -7578
# Row 132
$ws.Cells.Item(132, 8).Value = 1743.3846  # H132: was 1451.5454
$ws.Cells.Item(132, 9).Value = 1743.3846  # I132: was 1451.5454
$ws.Cells.Item(132, 11).Value = 5230.1538  # K132: was 4354.6362
$ws.Cells.Item(132, 13).Value = -2700.1538  # M132: was -1824.6362

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Cells.Item(20, 8).Value = 11749.5  # H20: was 13666.667
$ws.Cells.Item(20, 9).Value = 11666.333  # I20: was 13000
$ws.Cells.Item(20, 10).Value = 11999  # J20: was 15000
$ws.Cells.Item(20, 11).Value = 11666.333  # K20: was 13000
$ws.Cells.Item(20, 12).Value = 11999  # L20: was 15000
$ws.Cells.Item(20, 13).Value = -11419.333  # M20: was -12753
$ws.Cells.Item(20, 14).Value = -12493  # N20: was -15494
# Row 94
$ws.Cells.Item(94, 8).Value = 792.8461  # H94: was 761.5357
$ws.Cells.Item(94, 9).Value = 792.8461  # I94: was 801.0769
$ws.Cells.Item(94, 10).Value = 0  # J94: was 247.5
$ws.Cells.Item(94, 11).Value = 792.8461  # K94: was 801.0769
$ws.Cells.Item(94, 12).Value = 0  # L94: was 247.5
$ws.Cells.Item(94, 13).Value = -341.8461  # M94: was -350.0769
$ws.Cells.Item(94, 14).ClearContents()  # N94: was -1149.5
# Row 99
$ws.Cells.Item(99, 8).Value = 2654.04  # H99: was 2830.652
$ws.Cells.Item(99, 9).Value = 2308.7896  # I99: was 2507.1177
$ws.Cells.Item(99, 11).Value = 2308.7896  # K99: was 2507.1177
$ws.Cells.Item(99, 13).Value = -810.7896000000001  # M99: was -1009.1177

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Cells.Item(7, 8).Value = 184.22223  # H7: was 200.78947
$ws.Cells.Item(7, 10).Value = 500  # J7: was 499.5
$ws.Cells.Item(7, 12).Value = 500  # L7: was 499.5
$ws.Cells.Item(7, 14).Value = -726  # N7: was -725.5
# Row 22
$ws.Cells.Item(22, 8).Value = 633.3333  # H22: was 389.4
$ws.Cells.Item(22, 9).Value = 300  # I22: was 115.666664
$ws.Cells.Item(22, 11).Value = 300  # K22: was 115.666664
$ws.Cells.Item(22, 13).Value = 50  # M22: was 234.333336
# Row 112
$ws.Cells.Item(112, 8).Value = 13900  # H112: was 13949.5
$ws.Cells.Item(112, 10).Value = 13900  # J112: was 13949.5
$ws.Cells.Item(112, 12).Value = 13900  # L112: was 13949.5
$ws.Cells.Item(112, 14).Value = -16854  # N112: was -16903.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 116
$ws.Cells.Item(116, 8).Value = 1493.5  # H116: was 1497
$ws.Cells.Item(116, 9).Value = 1691.6666  # I116: was 1696.3334
$ws.Cells.Item(116, 11).Value = 5074.9998  # K116: was 5089.0002
$ws.Cells.Item(116, 13).Value = -1632.9998  # M116: was -1647.0002
# Row 129
$ws.Cells.Item(129, 8).Value = 2271.6667  # H129: was 2016.6923
$ws.Cells.Item(129, 9).Value = 1000  # I129: was 997.5
$ws.Cells.Item(129, 10).Value = 2526  # J129: was 2469.6667
$ws.Cells.Item(129, 11).Value = 3000  # K129: was 2992.5
$ws.Cells.Item(129, 12).Value = 7578  # L129: was 7409.000100000001
$ws.Cells.Item(129, 13).Value = 2000  # M129: was 2007.5
$ws.Cells.Item(129, 14).Value = -17578  # N129: was -17409.0001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Cells.Item(102, 8).Value = 0  # H102: was 310.8
$ws.Cells.Item(102, 9).Value = 0  # I102: was 310.8
$ws.Cells.Item(102, 11).Value = 0  # K102: was 310.8
$ws.Cells.Item(102, 13).ClearContents()  # M102: was 1311.2
# Row 107
$ws.Cells.Item(107, 8).Value = 1021.4  # H107: was 1037.1818
$ws.Cells.Item(107, 9).Value = 1086.1  # I107: was 1203.2222
$ws.Cells.Item(107, 10).Value = 892  # J107: was 922.2308
$ws.Cells.Item(107, 11).Value = 1086.1  # K107: was 1203.2222
$ws.Cells.Item(107, 12).Value = 892  # L107: was 922.2308
$ws.Cells.Item(107, 13).Value = 833.9000000000001  # M107: was 716.7778000000001
$ws.Cells.Item(107, 14).Value = -4732  # N107: was -4762.2308
# Row 113
$ws.Cells.Item(113, 8).Value = 4348.8125  # H113: was 5236.273
$ws.Cells.Item(113, 10).Value = 4145.769  # J113: was 5239.125
$ws.Cells.Item(113, 12).Value = 4145.769  # L113: was 5239.125
$ws.Cells.Item(113, 14).Value = -8485.769  # N113: was -9579.125
# Row 132
$ws.Cells.Item(132, 8).Value = 5719.6  # H132: was 3149.4443
$ws.Cells.Item(132, 9).Value = 1912  # I132: was 1779.5714
$ws.Cells.Item(132, 10).Value = 8258  # J132: was 7944
$ws.Cells.Item(132, 11).Value = 5736  # K132: was 5338.7142
$ws.Cells.Item(132, 12).Value = 24774  # L132: was 23832
$ws.Cells.Item(132, 13).Value = -3206  # M132: was -2808.7142
$ws.Cells.Item(132, 14).Value = -29834  # N132: was -28892

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 3
$ws.Cells.Item(3, 8).Value = 512  # H3: was 346.33334
$ws.Cells.Item(3, 9).Value = 24  # I3: was 19.5
$ws.Cells.Item(3, 11).Value = 24  # K3: was 19.5
$ws.Cells.Item(3, 13).Value = 88  # M3: was 92.5
# Row 7
$ws.Cells.Item(7, 8).Value = 2584.5  # H7: was 2621.4
$ws.Cells.Item(7, 9).Value = 2501.4  # I7: was 2376.75
$ws.Cells.Item(7, 10).Value = 3000  # J7: was 3600
$ws.Cells.Item(7, 11).Value = 2501.4  # K7: was 2376.75
$ws.Cells.Item(7, 12).Value = 3000  # L7: was 3600
$ws.Cells.Item(7, 13).Value = -2389.4  # M7: was -2264.75
$ws.Cells.Item(7, 14).Value = -3224  # N7: was -3824
# Row 15
$ws.Cells.Item(15, 8).Value = 512  # H15: was 346.33334
$ws.Cells.Item(15, 9).Value = 24  # I15: was 19.5
$ws.Cells.Item(15, 11).Value = 24  # K15: was 19.5
$ws.Cells.Item(15, 13).Value = 146  # M15: was 150.5
# Row 40
$ws.Cells.Item(40, 8).Value = 2683.3333  # H40: was 2236.1428
$ws.Cells.Item(40, 9).Value = 2683.3333  # I40: was 2236.1428
$ws.Cells.Item(40, 11).Value = 2683.3333  # K40: was 2236.1428
$ws.Cells.Item(40, 13).Value = -2547.3333  # M40: was -2100.1428
# Row 46
$ws.Cells.Item(46, 8).Value = 3119.25  # H46: was 3026.5908
$ws.Cells.Item(46, 9).Value = 2333  # I46: was 2499.75
$ws.Cells.Item(46, 10).Value = 3258  # J46: was 3143.6667
$ws.Cells.Item(46, 11).Value = 2333  # K46: was 2499.75
$ws.Cells.Item(46, 12).Value = 3258  # L46: was 3143.6667
$ws.Cells.Item(46, 13).Value = -2145  # M46: was -2311.75
$ws.Cells.Item(46, 14).Value = -3634  # N46: was -3519.6667
# Row 61
$ws.Cells.Item(61, 8).Value = 2565.6428  # H61: was 2825.9092
$ws.Cells.Item(61, 9).Value = 2637.182  # I61: was 3021.875
$ws.Cells.Item(61, 11).Value = 2637.182  # K61: was 3021.875
$ws.Cells.Item(61, 13).Value = -2435.182  # M61: was -2819.875
# Row 99
$ws.Cells.Item(99, 8).Value = 21923.334  # H99: was 21925.334
$ws.Cells.Item(99, 9).Value = 21923.334  # I99: was 21925.334
$ws.Cells.Item(99, 11).Value = 21923.334  # K99: was 21925.334
$ws.Cells.Item(99, 13).Value = -18928.334  # M99: was -18930.334
# Row 113
$ws.Cells.Item(113, 8).Value = 2565.6428  # H113: was 2825.9092
$ws.Cells.Item(113, 9).Value = 2637.182  # I113: was 3021.875
$ws.Cells.Item(113, 11).Value = 2637.182  # K113: was 3021.875
$ws.Cells.Item(113, 13).Value = -467.1819999999998  # M113: was -851.875
# Row 126
$ws.Cells.Item(126, 8).Value = 2584.5  # H126: was 2621.4
$ws.Cells.Item(126, 9).Value = 2501.4  # I126: was 2376.75
$ws.Cells.Item(126, 10).Value = 3000  # J126: was 3600
$ws.Cells.Item(126, 11).Value = 7504.200000000001  # K126: was 7130.25
$ws.Cells.Item(126, 12).Value = 9000  # L126: was 10800
$ws.Cells.Item(126, 13).Value = -5034.200000000001  # M126: was -4660.25
$ws.Cells.Item(126, 14).Value = -13940  # N126: was -15740

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Cells.Item(96, 8).Value = 1416.3334  # H96: was 1391.1818
$ws.Cells.Item(96, 9).Value = 1396.25  # I96: was 1361.25
$ws.Cells.Item(96, 10).Value = 1432.4  # J96: was 1408.2858
$ws.Cells.Item(96, 11).Value = 1396.25  # K96: was 1361.25
$ws.Cells.Item(96, 12).Value = 1432.4  # L96: was 1408.2858
$ws.Cells.Item(96, 13).Value = -23.25  # M96: was 11.75
$ws.Cells.Item(96, 14).Value = -4178.4  # N96: was -4154.2858
# Row 107
$ws.Cells.Item(107, 8).Value = 552.8  # H107: was 578.7143
$ws.Cells.Item(107, 9).Value = 423.8  # I107: was 482.25
$ws.Cells.Item(107, 11).Value = 1271.4  # K107: was 1446.75
$ws.Cells.Item(107, 13).Value = 648.5999999999999  # M107: was 473.25
# Row 119
$ws.Cells.Item(119, 8).Value = 0  # H119: was 11699
$ws.Cells.Item(119, 10).Value = 0  # J119: was 11699
$ws.Cells.Item(119, 12).Value = 0  # L119: was 11699
$ws.Cells.Item(119, 14).ClearContents()  # N119: was -21375
# Row 126
$ws.Cells.Item(126, 8).Value = 1869.8462  # H126: was 2025.3334
$ws.Cells.Item(126, 9).Value = 1230.3  # I126: was 1366.5555
$ws.Cells.Item(126, 11).Value = 3690.9  # K126: was 4099.666499999999
$ws.Cells.Item(126, 13).Value = -1220.9  # M126: was -1629.666499999999

Write-Host "Seraph profit sheets updated."
